# Insert 4 new price rows (one trading day's worth of Repollo quotes) just
# before the existing row 150, shifting rows 150:262 down to 154:266.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A150:A153").EntireRow.Insert()

# Populate the 4 freshly inserted rows with the new day's data.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Categoría ID,
#          G Categoría, H Variedad, I Calidad, J Volumen, K Precio mínimo,
#          L Precio máximo, M Precio promedio ponderado, N Unidad de
#          comercialización, O Origen, P Precio $/Kg, Q Kg o Unidades,
#          R Clasificación.

$rowsData = @(
    @{ Row=150; H="Copenhague"; I="Primera"; J=600;  K=900; L=1000; M=950 },
    @{ Row=151; H="Copenhague"; I="Segunda"; J=300;  K=800; L=800;  M=800 },
    @{ Row=152; H="Morada(o)";  I="Primera"; J=800;  K=900; L=1000; M=950 },
    @{ Row=153; H="Morada(o)";  I="Segunda"; J=400;  K=800; L=800;  M=800 }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Range("A$r").Value = 11
    $ws.Range("B$r").Value = "Vega Monumental Concepción"
    $ws.Range("C$r").Value = "Bíobío"
    $ws.Range("D$r").Value = 44574
    $ws.Range("E$r").Value = 8
    $ws.Range("F$r").Value = 100112006
    $ws.Range("G$r").Value = "Repollo"
    $ws.Range("H$r").Value = $rd.H
    $ws.Range("I$r").Value = $rd.I
    $ws.Range("J$r").Value = $rd.J
    $ws.Range("K$r").Value = $rd.K
    $ws.Range("L$r").Value = $rd.L
    $ws.Range("M$r").Value = $rd.M
    $ws.Range("N$r").Value = '$/unidad'
    $ws.Range("O$r").Value = "Región Metropolitana"
    $ws.Range("P$r").Value = $rd.M
    $ws.Range("Q$r").Value = 1
    $ws.Range("R$r").Value = "Hortaliza"
}
